$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.04110066666666667
$ws.Range("H2").Value = 0.123302
$ws.Range("I2").Value = 0.02671259512010182
$ws.Range("J2").Value = 0.02671259512010182
$ws.Range("M2").Value = 49.88013833333334
$ws.Range("N2").Value = 149.640415
$ws.Range("O2").Value = 0.9202778328538029
$ws.Range("P2").Value = 0.9202778328538028
$ws.Range("Q2").Value = 2.050106938925556
$ws.Range("R2").Value = 18.45096245033
$ws.Range("S2").Value = 0.02458300914702838
$ws.Range("T2").Value = 0.02458300914702838

$ws.Range("G3").Value = 0.04110066666666667
$ws.Range("H3").Value = 0.123302
$ws.Range("I3").Value = 0.02671259512010182
$ws.Range("J3").Value = 0.02671259512010182
$ws.Range("M3").Value = 0.6555886666666667
$ws.Range("O3").Value = 0.0120954700119653
$ws.Range("P3").Value = 0.0120954700119653
$ws.Range("Q3").Value = 0.02694513125911111
$ws.Range("R3").Value = 0.242506181332
$ws.Range("S3").Value = 0.0003231013932169622
$ws.Range("T3").Value = 0.0003231013932169622

$ws.Range("G4").Value = 0.04110066666666667
$ws.Range("H4").Value = 0.123302
$ws.Range("I4").Value = 0.02671259512010182
$ws.Range("J4").Value = 0.02671259512010182
$ws.Range("M4").Value = 0.126767
$ws.Range("N4").Value = 0.380301
$ws.Range("O4").Value = 0.002338823907379126
$ws.Range("P4").Value = 0.002338823907379127
$ws.Range("Q4").Value = 0.005210208211333333
$ws.Range("R4").Value = 0.046891873902
$ws.Range("S4").Value = 0.00006247605609503312
$ws.Range("T4").Value = 0.00006247605609503314

$ws.Range("G5").Value = 0.04110066666666667
$ws.Range("H5").Value = 0.123302
$ws.Range("I5").Value = 0.02671259512010182
$ws.Range("J5").Value = 0.02671259512010182
$ws.Range("M5").Value = 0.4520693333333334
$ws.Range("N5").Value = 1.356208
$ws.Range("O5").Value = 0.008340582048900294
$ws.Range("P5").Value = 0.008340582048900293
$ws.Range("Q5").Value = 0.01858035097955556
$ws.Range("R5").Value = 0.167223158816
$ws.Range("S5").Value = 0.0002227985913382629
$ws.Range("T5").Value = 0.0002227985913382628

$ws.Range("G6").Value = 0.04110066666666667
$ws.Range("H6").Value = 0.123302
$ws.Range("I6").Value = 0.02671259512010182
$ws.Range("J6").Value = 0.02671259512010182
$ws.Range("M6").Value = 2.307450666666667
$ws.Range("N6").Value = 6.922352
$ws.Range("O6").Value = 0.04257196892170599
$ws.Range("P6").Value = 0.04257196892170599
$ws.Range("Q6").Value = 0.09483776070044445
$ws.Range("R6").Value = 0.8535398463039999
$ws.Range("S6").Value = 0.00113720776927109
$ws.Range("T6").Value = 0.00113720776927109

$ws.Range("G7").Value = 0.04110066666666667
$ws.Range("H7").Value = 0.123302
$ws.Range("I7").Value = 0.02671259512010182
$ws.Range("J7").Value = 0.02671259512010182
$ws.Range("M7").Value = 0.7791593333333333
$ws.Range("N7").Value = 2.337478
$ws.Range("O7").Value = 0.01437532225624636
$ws.Range("P7").Value = 0.01437532225624636
$ws.Range("Q7").Value = 0.03202396803955555
$ws.Range("R7").Value = 0.288215712356
$ws.Range("S7").Value = 0.0003840021631520976
$ws.Range("T7").Value = 0.0003840021631520976

$ws.Range("I8").Value = 0.6362555311831452
$ws.Range("J8").Value = 0.636255531183145
$ws.Range("M8").Value = 49.88013833333334
$ws.Range("N8").Value = 149.640415
$ws.Range("O8").Value = 0.9202778328538029
$ws.Range("P8").Value = 0.9202778328538028
$ws.Range("Q8").Value = 48.83059371594889
$ws.Range("R8").Value = 439.4753434435401
$ws.Range("S8").Value = 0.58553186137847
$ws.Range("T8").Value = 0.5855318613784699

$ws.Range("I9").Value = 0.6362555311831452
$ws.Range("J9").Value = 0.636255531183145
$ws.Range("M9").Value = 0.6555886666666667
$ws.Range("O9").Value = 0.0120954700119653
$ws.Range("P9").Value = 0.0120954700119653
$ws.Range("Q9").Value = 0.6417942070017778
$ws.Range("R9").Value = 5.776147863016
$ws.Range("S9").Value = 0.007695809697372784
$ws.Range("T9").Value = 0.007695809697372783

$ws.Range("I10").Value = 0.6362555311831452
$ws.Range("J10").Value = 0.636255531183145
$ws.Range("M10").Value = 0.126767
$ws.Range("N10").Value = 0.380301
$ws.Range("O10").Value = 0.002338823907379126
$ws.Range("P10").Value = 0.002338823907379127
$ws.Range("Q10").Value = 0.1240996532973333
$ws.Range("R10").Value = 1.116896879676
$ws.Range("S10").Value = 0.001488089647533345
$ws.Range("T10").Value = 0.001488089647533345

$ws.Range("I11").Value = 0.6362555311831452
$ws.Range("J11").Value = 0.636255531183145
$ws.Range("M11").Value = 0.4520693333333334
$ws.Range("N11").Value = 1.356208
$ws.Range("O11").Value = 0.008340582048900294
$ws.Range("P11").Value = 0.008340582048900293
$ws.Range("Q11").Value = 0.4425571918008889
$ws.Range("R11").Value = 3.983014726208
$ws.Range("S11").Value = 0.005306741461899662
$ws.Range("T11").Value = 0.00530674146189966

$ws.Range("I12").Value = 0.6362555311831452
$ws.Range("J12").Value = 0.636255531183145
$ws.Range("M12").Value = 2.307450666666667
$ws.Range("N12").Value = 6.922352
$ws.Range("O12").Value = 0.04257196892170599
$ws.Range("P12").Value = 0.04257196892170599
$ws.Range("Q12").Value = 2.258898828039111
$ws.Range("R12").Value = 20.330089452352
$ws.Range("S12").Value = 0.02708665069979239
$ws.Range("T12").Value = 0.02708665069979239

$ws.Range("I13").Value = 0.6362555311831452
$ws.Range("J13").Value = 0.636255531183145
$ws.Range("M13").Value = 0.7791593333333333
$ws.Range("N13").Value = 2.337478
$ws.Range("O13").Value = 0.01437532225624636
$ws.Range("P13").Value = 0.01437532225624636
$ws.Range("Q13").Value = 0.7627647820808888
$ws.Range("R13").Value = 6.864883038727999
$ws.Range("S13").Value = 0.009146378298076914
$ws.Range("T13").Value = 0.009146378298076912

$ws.Range("G14").Value = 0.5185656666666667
$ws.Range("H14").Value = 1.555697
$ws.Range("I14").Value = 0.3370318736967531
$ws.Range("J14").Value = 0.3370318736967531
$ws.Range("M14").Value = 49.88013833333334
$ws.Range("N14").Value = 149.640415
$ws.Range("O14").Value = 0.9202778328538029
$ws.Range("P14").Value = 0.9202778328538028
$ws.Range("Q14").Value = 25.86612718825056
$ws.Range("R14").Value = 232.795144694255
$ws.Range("S14").Value = 0.3101629623283045
$ws.Range("T14").Value = 0.3101629623283045

$ws.Range("G15").Value = 0.5185656666666667
$ws.Range("H15").Value = 1.555697
$ws.Range("I15").Value = 0.3370318736967531
$ws.Range("J15").Value = 0.3370318736967531
$ws.Range("M15").Value = 0.6555886666666667
$ws.Range("O15").Value = 0.0120954700119653
$ws.Range("P15").Value = 0.0120954700119653
$ws.Range("Q15").Value = 0.3399657739891112
$ws.Range("R15").Value = 3.059691965902001
$ws.Range("S15").Value = 0.004076558921375553
$ws.Range("T15").Value = 0.004076558921375553

$ws.Range("G16").Value = 0.5185656666666667
$ws.Range("H16").Value = 1.555697
$ws.Range("I16").Value = 0.3370318736967531
$ws.Range("J16").Value = 0.3370318736967531
$ws.Range("M16").Value = 0.126767
$ws.Range("N16").Value = 0.380301
$ws.Range("O16").Value = 0.002338823907379126
$ws.Range("P16").Value = 0.002338823907379127
$ws.Range("Q16").Value = 0.06573701386633334
$ws.Range("R16").Value = 0.591633124797
$ws.Range("S16").Value = 0.0007882582037507482
$ws.Range("T16").Value = 0.0007882582037507483

$ws.Range("G17").Value = 0.5185656666666667
$ws.Range("H17").Value = 1.555697
$ws.Range("I17").Value = 0.3370318736967531
$ws.Range("J17").Value = 0.3370318736967531
$ws.Range("M17").Value = 0.4520693333333334
$ws.Range("N17").Value = 1.356208
$ws.Range("O17").Value = 0.008340582048900294
$ws.Range("P17").Value = 0.008340582048900293
$ws.Range("Q17").Value = 0.2344276352195556
$ws.Range("R17").Value = 2.109848716976
$ws.Range("S17").Value = 0.00281104199566237
$ws.Range("T17").Value = 0.00281104199566237

$ws.Range("G18").Value = 0.5185656666666667
$ws.Range("H18").Value = 1.555697
$ws.Range("I18").Value = 0.3370318736967531
$ws.Range("J18").Value = 0.3370318736967531
$ws.Range("M18").Value = 2.307450666666667
$ws.Range("N18").Value = 6.922352
$ws.Range("O18").Value = 0.04257196892170599
$ws.Range("P18").Value = 0.04257196892170599
$ws.Range("Q18").Value = 1.196564693260445
$ws.Range("R18").Value = 10.769082239344
$ws.Range("S18").Value = 0.01434811045264251
$ws.Range("T18").Value = 0.01434811045264251

$ws.Range("G19").Value = 0.5185656666666667
$ws.Range("H19").Value = 1.555697
$ws.Range("I19").Value = 0.3370318736967531
$ws.Range("J19").Value = 0.3370318736967531
$ws.Range("M19").Value = 0.7791593333333333
$ws.Range("N19").Value = 2.337478
$ws.Range("O19").Value = 0.01437532225624636
$ws.Range("P19").Value = 0.01437532225624636
$ws.Range("Q19").Value = 0.4040452791295556
$ws.Range("R19").Value = 3.636407512166
$ws.Range("S19").Value = 0.004844941795017346
$ws.Range("T19").Value = 0.004844941795017346
